# ReadFromFile: populate the "user_details" sheet with the test data used by
# the Checkout "Your Information" page (First Name / Last Name / Zip-Postal
# Code), styled like the existing "login" sheet, and leave that sheet active
# since it's the one the author was last working on.

$wb = $excel.ActiveWorkbook

$loginSheet   = $wb.Worksheets.Item("login")
$detailsSheet = $wb.Worksheets.Item("user_details")

# Data row first (matches the order the values were typed in originally),
# then the header row.
$detailsSheet.Range("A2").Value = "Tami"
$detailsSheet.Range("B2").Value = "Tam"

$detailsSheet.Range("A1").Value = "First Name"
$detailsSheet.Range("B1").Value = "Last Name"
$detailsSheet.Range("C1").Value = "Zip/Postal Code"

$detailsSheet.Range("C2").Value = 2121

# Give the new cells the same bordered look used for the "login" sheet's
# data rows.
$loginSheet.Range("A2:B2").Copy()
$detailsSheet.Range("A1:C2").PasteSpecial(-4122) # xlPasteFormats

# Switch focus to "user_details" - it becomes the active tab/sheet, and
# leave the selection where the author last clicked.
$detailsSheet.Activate()
$detailsSheet.Range("D20").Select()

$wb.Save()
